# Update cryptocurrency price/volume data per upstream source refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.848.14'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '2.355.50'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''504.00'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '''129.93'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '2.369.63'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = '2.774.64'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '55.792.23'
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('D16').Value = '''21.44'
$ws.Range('E16').Value = '  -1.22%  '
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '2.371.90'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').Value = '''310.11'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = '''6.21'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').Value = '''65.28'
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  -3.14%  '
$ws.Range('D28').Value = '''7.12'
$ws.Range('E28').Value = '  -3.34%  '
$ws.Range('D29').Value = '''170.76'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').Value = '0.0₃0705'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '''0.996'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('E35').Value = '  -4.97%  '
$ws.Range('D36').Value = '''17.64'
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('D38').Value = '''0.836'
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('D39').Value = '''3.64'
$ws.Range('E39').Value = '  -4.37%  '
$ws.Range('D40').Value = '''36.17'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('E41').Value = '  -2.91%  '
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').Value = '''125.77'
$ws.Range('E44').Value = '  -5.05%  '
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '''239.46'
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').Value = '''16.85'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('D51').Value = '''16.54'
$ws.Range('E51').Value = '  -3.55%  '
